# Auto-generated: apply scheduled market-data refresh to profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 718.7778
$ws.Range("I12").Value = 679.3333
$ws.Range("J12").Value = 797.6667
$ws.Range("K12").Value = 679.3333
$ws.Range("L12").Value = 797.6667
$ws.Range("M12").Value = -509.3333
$ws.Range("N12").Value = -1137.6667

# Row 32
$ws.Range("H32").Value = 3783.3572
$ws.Range("I32").Value = 2697.8
$ws.Range("J32").Value = 4019.348
$ws.Range("K32").Value = 2697.8
$ws.Range("L32").Value = 4019.348
$ws.Range("M32").Value = -2371.8
$ws.Range("N32").Value = -4671.348

# Row 86
$ws.Range("H86").Value = 4090.4062
$ws.Range("I86").Value = 3652.1904
$ws.Range("J86").Value = 4927
$ws.Range("K86").Value = 3652.1904
$ws.Range("L86").Value = 4927
$ws.Range("M86").Value = -2529.1904
$ws.Range("N86").Value = -7173

# Row 89
$ws.Range("H89").Value = 4090.4062
$ws.Range("I89").Value = 3652.1904
$ws.Range("J89").Value = 4927
$ws.Range("K89").Value = 18260.952
$ws.Range("L89").Value = 24635
$ws.Range("M89").Value = -12644.952
$ws.Range("N89").Value = -35867

# Row 138
$ws.Range("H138").Value = 8766.166999999999
$ws.Range("I138").Value = 6498
$ws.Range("J138").Value = 9049.6875
$ws.Range("K138").Value = 19494
$ws.Range("L138").Value = 27149.0625
$ws.Range("M138").Value = -14354
$ws.Range("N138").Value = -37429.0625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17256.857
$ws.Range("I32").Value = 14064.328
$ws.Range("J32").Value = 32973.92
$ws.Range("K32").Value = 14064.328
$ws.Range("L32").Value = 32973.92
$ws.Range("M32").Value = -13777.328
$ws.Range("N32").Value = -33547.92

# Row 61
$ws.Range("H61").Value = 17247228
$ws.Range("I61").Value = 22733120
$ws.Range("K61").Value = 22733120
$ws.Range("M61").Value = -22732908

# Row 74
$ws.Range("H74").Value = 63298.535
$ws.Range("I74").Value = 4724.4165
$ws.Range("J74").Value = 297595
$ws.Range("K74").Value = 4724.4165
$ws.Range("L74").Value = 297595
$ws.Range("M74").Value = -3850.4165
$ws.Range("N74").Value = -299343

# Row 77
$ws.Range("H77").Value = 63298.535
$ws.Range("I77").Value = 4724.4165
$ws.Range("J77").Value = 297595
$ws.Range("K77").Value = 23622.0825
$ws.Range("L77").Value = 1487975
$ws.Range("M77").Value = -19254.0825
$ws.Range("N77").Value = -1496711

# Row 96
$ws.Range("H96").Value = 27000
$ws.Range("J96").Value = 27000
$ws.Range("L96").Value = 27000
$ws.Range("N96").Value = -32492

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 132
$ws.Range("H132").Value = 37032.1
$ws.Range("I132").Value = 2788.8572
$ws.Range("J132").Value = 68992.47
$ws.Range("K132").Value = 8366.571599999999
$ws.Range("L132").Value = 206977.41
$ws.Range("M132").Value = -5836.571599999999
$ws.Range("N132").Value = -212037.41

# Row 136
$ws.Range("H136").Value = 17247228
$ws.Range("I136").Value = 22733120
$ws.Range("K136").Value = 68199360
$ws.Range("M136").Value = -68196810

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 9301.091
$ws.Range("I58").Value = 11090.059
$ws.Range("K58").Value = 11090.059
$ws.Range("M58").Value = -10887.059

# Row 132
$ws.Range("H132").Value = 120332.69
$ws.Range("I132").Value = 73674.28999999999
$ws.Range("K132").Value = 221022.87
$ws.Range("M132").Value = -218492.87

# Row 134
$ws.Range("H134").Value = 2626.4075
$ws.Range("I134").Value = 1871.36
$ws.Range("J134").Value = 3277.3103
$ws.Range("K134").Value = 5614.08
$ws.Range("L134").Value = 9831.930899999999
$ws.Range("M134").Value = -3079.08
$ws.Range("N134").Value = -14901.9309

# Row 136
$ws.Range("H136").Value = 9301.091
$ws.Range("I136").Value = 11090.059
$ws.Range("K136").Value = 33270.177
$ws.Range("M136").Value = -30720.177

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 57063.277
$ws.Range("J5").Value = 101995.7
$ws.Range("L5").Value = 305987.1
$ws.Range("N5").Value = -306211.1

# Row 56
$ws.Range("H56").Value = 22733672
$ws.Range("I56").Value = 22733672
$ws.Range("K56").Value = 22733672
$ws.Range("M56").Value = -22733142

# Row 122
$ws.Range("H122").Value = 2238.7778
$ws.Range("J122").Value = 3312.5
$ws.Range("L122").Value = 29812.5
$ws.Range("N122").Value = -34712.5

# Row 135
$ws.Range("H135").Value = 57063.277
$ws.Range("J135").Value = 101995.7
$ws.Range("L135").Value = 917961.2999999999
$ws.Range("N135").Value = -923031.2999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 16680092
$ws.Range("I70").Value = 20004470
$ws.Range("K70").Value = 20004470
$ws.Range("M70").Value = -20004200

# Row 73
$ws.Range("H73").Value = 16680092
$ws.Range("I73").Value = 20004470
$ws.Range("K73").Value = 20004470
$ws.Range("M73").Value = -20003534

# Row 80
$ws.Range("H80").Value = 2241890.5
$ws.Range("I80").Value = 3766493.5
$ws.Range("J80").Value = 526712.4
$ws.Range("K80").Value = 3766493.5
$ws.Range("L80").Value = 526712.4
$ws.Range("M80").Value = -3765495.5
$ws.Range("N80").Value = -528708.4

# Row 83
$ws.Range("H83").Value = 2241890.5
$ws.Range("I83").Value = 3766493.5
$ws.Range("J83").Value = 526712.4
$ws.Range("K83").Value = 18832467.5
$ws.Range("L83").Value = 2633562
$ws.Range("M83").Value = -18827475.5
$ws.Range("N83").Value = -2643546

# Row 132
$ws.Range("H132").Value = 4989.75
$ws.Range("I132").Value = 4989.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14969.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12439.25
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 11394.3
$ws.Range("I7").Value = 10986.75
$ws.Range("J7").Value = 11666
$ws.Range("K7").Value = 10986.75
$ws.Range("L7").Value = 11666
$ws.Range("M7").Value = -10874.75
$ws.Range("N7").Value = -11890

# Row 55
$ws.Range("H55").Value = 1453.8462
$ws.Range("I55").Value = 1388.1177
$ws.Range("K55").Value = 1388.1177
$ws.Range("M55").Value = -1215.1177

# Row 68
$ws.Range("H68").Value = 4724.5
$ws.Range("I68").Value = 2479.2
$ws.Range("K68").Value = 2479.2
$ws.Range("M68").Value = -1730.2

# Row 71
$ws.Range("H71").Value = 4724.5
$ws.Range("I71").Value = 2479.2
$ws.Range("K71").Value = 12396
$ws.Range("M71").Value = -8652

# Row 126
$ws.Range("H126").Value = 11394.3
$ws.Range("I126").Value = 10986.75
$ws.Range("J126").Value = 11666
$ws.Range("K126").Value = 32960.25
$ws.Range("L126").Value = 34998
$ws.Range("M126").Value = -30490.25
$ws.Range("N126").Value = -39938

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2316.7778
$ws.Range("I96").Value = 1586.5
$ws.Range("J96").Value = 2525.4285
$ws.Range("K96").Value = 1586.5
$ws.Range("L96").Value = 2525.4285
$ws.Range("M96").Value = -213.5
$ws.Range("N96").Value = -5271.4285

# Row 109
$ws.Range("H109").Value = 45377.92
$ws.Range("J109").Value = 45377.92
$ws.Range("L109").Value = 45377.92
$ws.Range("N109").Value = -48151.92

# Row 115
$ws.Range("H115").Value = 49071.54
$ws.Range("J115").Value = 49071.54
$ws.Range("L115").Value = 49071.54
$ws.Range("N115").Value = -52205.54

# Row 132
$ws.Range("H132").Value = 36109520
$ws.Range("I132").Value = 50007372
$ws.Range("J132").Value = 1364889
$ws.Range("K132").Value = 150022116
$ws.Range("L132").Value = 4094667
$ws.Range("M132").Value = -150019586
$ws.Range("N132").Value = -4099727

